# Insert a new "Industry" column as column C, shifting the existing
# Mutual Fund / Status / Jan_2026 / Dec_2025 / Oct_2025 / MoM / QoQ
# columns (old C:I) one position to the right (new D:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C. This shifts columns C.. to D..
# and carries the existing cell formatting (e.g. the header style) along
# with the shifted cells, while the newly inserted column starts blank.
$ws.Columns.Item(3).Insert()

# Header for the new column, matching the style of the other header cells.
$ws.Cells.Item(1, 3).Value = "Industry"
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 3).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Industry values for each holding row (row number -> industry name).
$industries = @{
    2  = "Banks"
    3  = "Banks"
    4  = "Construction"
    5  = "Transport Services"
    6  = "Retailing"
    7  = "IT - Software"
    8  = "Finance"
    9  = "Aerospace & Defense"
    10 = "Automobiles"
    11 = "Electrical Equipment"
    12 = "Capital Markets"
    13 = "Consumer Durables"
    14 = "IT - Software"
    15 = "Transport Services"
    16 = "Healthcare Services"
    17 = "Beverages"
    18 = "Automobiles"
    19 = "Agricultural, Commercial & Construction Vehicles"
    20 = "Paper, Forest & Jute Products"
    21 = "Pharmaceuticals & Biotechnology"
    22 = "Automobiles"
    23 = "Agricultural, Commercial & Construction Vehicles"
    24 = "Capital Markets"
    25 = "Pharmaceuticals & Biotechnology"
    26 = "Retailing"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
